$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-12 (text + value) - these already carry the
# correct "s=1" style, so a plain Value assignment is sufficient.
$ws.Range("A2").Value = "IX.  Doenças do aparelho circulatório"
$ws.Range("B2").Value = 518254
$ws.Range("A3").Value = "II.  Neoplasias (tumores)"
$ws.Range("B3").Value = 358364
$ws.Range("A4").Value = "X.   Doenças do aparelho respiratório"
$ws.Range("B4").Value = 204420
$ws.Range("A5").Value = "XX.  Causas externas de morbidade e mortalidade"
$ws.Range("B5").Value = 143006
$ws.Range("A6").Value = "IV.  Doenças endócrinas nutricionais e metabólicas"
$ws.Range("B6").Value = 93801
$ws.Range("A7").Value = "I.   Algumas doenças infecciosas e parasitárias"
$ws.Range("B7").Value = 92615
$ws.Range("A8").Value = "XVIII.Sint sinais e achad anorm ex clín e laborat"
$ws.Range("B8").Value = 83448
$ws.Range("A9").Value = "XI.  Doenças do aparelho digestivo"
$ws.Range("B9").Value = 80315
$ws.Range("A10").Value = "VI.  Doenças do sistema nervoso"
$ws.Range("B10").Value = 49264
$ws.Range("A11").Value = "XIV. Doenças do aparelho geniturinário"
$ws.Range("B11").Value = 34486
$ws.Range("A12").Value = "V.   Transtornos mentais e comportamentais"
$ws.Range("B12").Value = 13382

# Copy the style of an existing formatted cell (A12) so the newly
# added rows (13-20) get the identical "s=1" formatting (bold font,
# border, centered/top alignment) instead of a brand-new style index.
$styleSource = $ws.Range("A12")

$ws.Range("A13").Value = "XIII.Doenças sist osteomuscular e tec conjuntivo"
$ws.Range("B13").Value = 6616
$styleSource.Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null

$ws.Range("A14").Value = "III. Doenças sangue órgãos hemat e transt imunitár"
$ws.Range("B14").Value = 6494
$styleSource.Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

$ws.Range("A15").Value = "XII. Doenças da pele e do tecido subcutâneo"
$ws.Range("B15").Value = 3661
$styleSource.Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null

$ws.Range("A16").Value = "XVII.Malf cong deformid e anomalias cromossômicas"
$ws.Range("B16").Value = 2408
$styleSource.Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

$ws.Range("A17").Value = "XV.  Gravidez parto e puerpério"
$ws.Range("B17").Value = 2063
$styleSource.Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null

$ws.Range("A18").Value = "XVI. Algumas afec originadas no período perinatal"
$ws.Range("B18").Value = 909
$styleSource.Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null

$ws.Range("A19").Value = "VIII.Doenças do ouvido e da apófise mastóide"
$ws.Range("B19").Value = 123
$styleSource.Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null

$ws.Range("A20").Value = "VII. Doenças do olho e anexos"
$ws.Range("B20").Value = 11
$styleSource.Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Refresh the sheet dimension to cover the new rows
$ws.Range("A1:B20").Select() | Out-Null